$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the previous day's two rows (4025 -> 四方坪站/高岭站)
# into the two new rows so the new cells pick up the existing number-format
# styles (date style, currency style, integer style) instead of creating new ones.
$ws.Range("A6:F7").Copy() | Out-Null
$ws.Range("A8:F9").PasteSpecial(-4122) | Out-Null

# Row 8: 2026-01-04 (serial 46026) - 四方坪站
$ws.Range("A8").Value = 46026
$ws.Range("B8").Value = "四方坪站"
$ws.Range("C8").Value = 13588.81
$ws.Range("D8").Value = 9759.66
$ws.Range("E8").Value = 3386.45
$ws.Range("F8").Value = 566

# Row 9: 2026-01-04 (serial 46026) - 高岭站
$ws.Range("A9").Value = 46026
$ws.Range("B9").Value = "高岭站"
$ws.Range("C9").Value = 6952.05
$ws.Range("D9").Value = 6137.13
$ws.Range("E9").Value = 1815.98
$ws.Range("F9").Value = 236

# Match the recorded selection state
$ws.Range("J15").Select() | Out-Null
